# Move RHPF to shift hydrogen to steam methane reforming; fix default
# hydrogen shifting FoPITY settings to 1 from 2025-2050; FoPITY edits for
# NDC scenario.

$wb = $excel.ActiveWorkbook

# --- "About" sheet: update description of the recipient pathway --------
$about = $wb.Worksheets.Item("About")
$about.Range("A12").Value = "steam methane reforming."

# --- "RHPF" sheet: shift the default recipient pathway ------------------
$rhpf = $wb.Worksheets.Item("RHPF")

# Previously every "from" pathway (columns B:H) defaulted its hydrogen
# fully (value 1) to row 7, "electrolysis with guaranteed clean
# electricity." Now it should default fully (value 1) to row 3, "natural
# gas reforming" (i.e. steam methane reforming), and row 7 drops to 0.
$rhpf.Range("B3:H3").Value = 1
$rhpf.Range("B7:H7").Value = 0

# Match formatting: row 3 (now populated like the rest of the matrix)
# picks up the right-aligned numeric style already used by columns B:F,
# while row 7 (now zero, same as the other non-populated rows) reverts
# to the default/general style used elsewhere in the table.
$rhpf.Range("G3:H3").HorizontalAlignment = -4152
$rhpf.Range("G7:H7").ClearFormats()

# Leave the selection on the row that was just edited, as the author did.
$rhpf.Activate() | Out-Null
$rhpf.Range("B7:H7").Select() | Out-Null
